$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data cleanup: trim stray leading/trailing spaces baked into the names ---
# ("Amir " -> "Amir", " Kumar" -> "Kumar", "Rishi " -> "Rishi") so the actor's
# name and the expected_response JSON agree on value, not just key.
$ws.Range("B2").Value = "Kumar"
$ws.Range("A3").Value = "Amir"
$ws.Range("A11").Value = "Rishi"

# --- Formatting cleanup for column C (expected_response) ---
# Column C no longer carries a blanket "wrap text" style; only the sample
# row (row 2) keeps it, sized so the wrapped JSON is fully visible.
$ws.Columns.Item(3).ClearFormats()

# Clearing the column format also wiped C1's header styling - restore it by
# copying the format from the (untouched) A1 header cell.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply wrap text to the example cell and grow its row to fit the JSON.
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 57.6

# --- Selection / view cleanup ---
$null = $ws.Range("B3").Select()
